# Fix a handful of typos found across the deck.
#   Webscrapping -> Webscraping   (slide 2)
#   Page Jaune    -> Pages Jaunes (slide 4)
#   Scorring      -> Scoring      (slide 5)
#
# Walk every shape's text runs on every slide and swap the exact text of
# any run that matches one of the known typos. Editing Run.Text (rather
# than the whole TextRange/shape) keeps each run's existing formatting
# (rPr) untouched, matching how PowerPoint itself records such a fix.

$p = $ppt.ActivePresentation

$replacements = @{
    "Webscrapping" = "Webscraping"
    "Page Jaune"   = "Pages Jaunes"
    "Scorring"     = "Scoring"
}

foreach ($s in $p.Slides) {
    foreach ($shape in $s.Shapes) {
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }

        foreach ($para in $tf.TextRange.Paragraphs()) {
            foreach ($run in $para.Runs()) {
                $newText = $replacements[$run.Text]
                if ($newText) {
                    $run.Text = $newText
                }
            }
        }
    }
}
